$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row:
#      "<Name>_old" -> "<Name>_FV2404"
#      "<Name>_new" -> "<Name>_FV2410"
#    Column K ("diff") stays the same.
# ---------------------------------------------------------------------------
$lastCol = 21 # columns A..U
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = [string]$cell.Value2
    if ($header -like "*_old") {
        $newHeader = $header.Substring(0, $header.Length - 4) + "_FV2404"
        $cell.Value = $newHeader
    } elseif ($header -like "*_new") {
        $newHeader = $header.Substring(0, $header.Length - 4) + "_FV2410"
        $cell.Value = $newHeader
    }
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (split after row 1, keep left pane at column A).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the data range into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$lastRow = $ws.UsedRange.Rows.Count
$tableRange = $ws.Range("A1").Resize($lastRow, $lastCol)
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
